$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = "Lenovo ThinkPad R61e"
$ws.Range("B4").Select()
